# Fruta / hortaliza, semanal
#
# Insert a new weekly data row at row 65 (pushing the existing rows 65-120
# down to 66-121, as a plain Excel row insert would), then populate the
# newly-inserted row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65; Excel shifts rows 65..120
# down to 66..121 and expands the used range automatically.
$ws.Rows.Item(65).Insert()

# Fill in the values for the freshly inserted row 65.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44512
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100112025
$ws.Range("J65").Value = "Frutilla"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 600
$ws.Range("N65").Value = 9000
$ws.Range("O65").Value = 10000
$ws.Range("P65").Value = 9500
$ws.Range("Q65").Value = '$/bandeja 7 kilos'
$ws.Range("R65").Value = "Provincia de Melipilla"
$ws.Range("S65").Value = 1357
$ws.Range("T65").Value = 7
